# ajout assets et credentials
#
# Adds two new Asset rows (SeLoger + helloPret URL credentials) to the
# "Assets" config sheet, and makes that sheet the active/selected tab
# (matching the author's last on-screen state when the workbook was saved).

$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# --- New Asset rows -------------------------------------------------------
# Row 2: seLoger URL asset
$wsAssets.Range("A2").Value = "url_seLoger"
$wsAssets.Range("B2").Value = "url_seLoger "
$wsAssets.Range("C2").Value = "logement"
$wsAssets.Range("D2").Value = "url de SeLoger"

# Row 3: helloPret URL asset
$wsAssets.Range("A3").Value = "url_helloPret "
$wsAssets.Range("B3").Value = "url_helloPret "
$wsAssets.Range("C3").Value = "logement"
$wsAssets.Range("D3").Value = "url de helloPret"

# Match the wrapped-text style already used for the descriptive columns
# (Name / Asset / Description) elsewhere in this workbook.
$wsAssets.Range("A2:B3").WrapText = $true
$wsAssets.Range("D2:D3").WrapText = $true

# --- Active tab / selection -------------------------------------------
# The workbook was left with the Assets sheet active and A4 selected.
$wsAssets.Activate() | Out-Null
$wsAssets.Range("A4").Select() | Out-Null

# --- Minor page setup normalisation (harmless, matches resave defaults) ---
foreach ($ws in @($wsSettings, $wsConstants, $wsAssets)) {
    $ws.PageSetup.FitToPagesWide = 1
    $ws.PageSetup.FitToPagesTall = 1
    $ws.PageSetup.HeaderMargin = 36.85
    $ws.PageSetup.FooterMargin = 36.85
}
